# "process updated for final maxico"
# Week 05 -> Week 07 rollover: update the MEXICO folder-path strings on the
# Settings sheet, bump the WorkingStartDate/WorkingEndDate on the Constants
# sheet, and leave the workbook with Constants as the last-active sheet
# (matching the saved selection/activeTab state captured in the diff).

$wb = $excel.ActiveWorkbook
$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------
# Settings sheet: "05 Envío/Envio Semana 05" -> "07 Envío/Envio Semana 07"
# ---------------------------------------------------------------------
$wsSettings.Activate()

$baseDatos    = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/MÉXICO/Base de Datos"
$consolidado  = "/Planeacion/0.Envios TS/2022/07 Envio Semana 07/MÉXICO/Base de Datos/Consolidado"
$vip          = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/MÉXICO/VIP"
$tradicional  = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/MÉXICO/Tradicional"
$miCine       = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/MÉXICO/Mi Cine"
$atmosfera    = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/MÉXICO/Atmosfera"
$exportadas   = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/MÉXICO/Base de Datos/Exportadas"

$wsSettings.Range("B2").Value  = $baseDatos
$wsSettings.Range("B6").Value  = $baseDatos
$wsSettings.Range("B20").Value = $baseDatos

$wsSettings.Range("B8").Value  = $consolidado

$wsSettings.Range("B10").Value = $vip
$wsSettings.Range("B16").Value = $vip

$wsSettings.Range("B12").Value = $tradicional
$wsSettings.Range("B14").Value = $miCine
$wsSettings.Range("B18").Value = $atmosfera

$wsSettings.Range("B23").Value = $exportadas
$wsSettings.Range("B26").Value = $exportadas
$wsSettings.Range("B29").Value = $exportadas
$wsSettings.Range("B33").Value = $exportadas

# Leaves this sheet's saved selection on B2 (was B5).
$wsSettings.Range("B2").Select()

# ---------------------------------------------------------------------
# Constants sheet: WorkingStartDate / WorkingEndDate -> next week
# ---------------------------------------------------------------------
$wsConstants.Activate()

$wsConstants.Range("B24").Value = 44599
$wsConstants.Range("B25").Value = 44605

# Constants becomes the last active sheet/tab, selection moves to B24
# (was B26), matching the saved workbook + sheet view state.
$wsConstants.Range("B24").Select()
